$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts old B/C/D -> C/D/E)
$ws.Columns.Item(2).Insert() | Out-Null

# New header for the inserted column
$ws.Range("B1").Value = "_requirements"

# Widen the new column to fit the longer requirement notes
$ws.Columns.Item(2).ColumnWidth = 44.62

# Rows with no requirement note: make sure no stray formatting remains
$ws.Range("B5").Clear() | Out-Null
$ws.Range("B11").Clear() | Out-Null

# Fill in the requirement notes on the relevant rows
$ws.Range("B6").Value = "l10n_it_reverse_charge"
$ws.Range("B7").Value = "l10n_it_reverse_charge"
$ws.Range("B8").Value = "l10n_it_split_payment"
$ws.Range("B9").Value = "l10n_it_dichiarazione_intento or l10n_it_lettera_intento "

# B6:B8 keep the plain (unstyled) look used elsewhere in the sheet
$ws.Range("D2").Copy() | Out-Null
$ws.Range("B6:B8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Match the workbook's final selection
$ws.Range("B10").Select() | Out-Null
